# "Fruta / hortaliza, semanal" — insert the latest weekly price record for
# Cilantro (Agrícola del Norte S.A. de Arica) as a new row 65, pushing the
# existing historical rows (old 65..89) down to (new 66..90).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 65, shifting rows 65:89 down to 66:90.
$ws.Rows.Item(65).Insert()

# Populate the new row 65 with the new weekly observation.
$ws.Range("A65").Value = 1
$ws.Range("B65").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C65").Value = "Arica y Parinacota"
$ws.Range("D65").Value = 44837
$ws.Range("E65").Value = 15
$ws.Range("F65").Value = 100112040
$ws.Range("G65").Value = "Cilantro"
$ws.Range("H65").Value = "Sin especificar"
$ws.Range("I65").Value = "Primera"
$ws.Range("J65").Value = 300
$ws.Range("K65").Value = 500
$ws.Range("L65").Value = 600
$ws.Range("M65").Value = 550
$ws.Range("N65").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O65").Value = "Región de Arica y Parinacota"
$ws.Range("P65").Value = 275
$ws.Range("Q65").Value = 2
$ws.Range("R65").Value = "Hortaliza"
